$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) from the last existing data row (row 8) down
# into the three new rows (9-11) so the new cells match the existing
# table's look (bold/bordered id column, date-formatted date column, etc.)
$ws.Range("A8:M8").Copy()
$ws.Range("A9:M9").PasteSpecial(-4122)
$ws.Range("A8:M8").Copy()
$ws.Range("A10:M10").PasteSpecial(-4122)
$ws.Range("A8:M8").Copy()
$ws.Range("A11:M11").PasteSpecial(-4122)

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "no"
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = "yes"
$ws.Range("G9").Value = 43690
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 30
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = "no"
$ws.Range("L9").Value = "exercise"
$ws.Range("M9").Value = 24

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "no"
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = "yes"
$ws.Range("G10").Value = 43740
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 15
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = "no"
$ws.Range("L10").Value = "work"
$ws.Range("M10").Value = 12

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "yes"
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = "yes"
$ws.Range("G11").Value = 43709
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 60
$ws.Range("J11").Value = 5
$ws.Range("K11").Value = "yes"
$ws.Range("L11").Value = "work"
$ws.Range("M11").Value = 12
